$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("sergio", "nacho", "nauce", "borja", "marcelo", "jorge")

$row = 27
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = 5000
    $row = $row + 1
}

$ws.Range("D33").Select()
